# Upload new version with timestamp
#
# Business changes captured by the diff:
#  - MINALAX 10 TABLETS row (row 12): current balance "2:0" -> "1:0",
#    sale price "18.0000" -> "36.0000", number of transactions "1:0" -> "2:0"
#  - Totals row (P15) recomputed to reflect the new sale price (+18)
#  - The "generated at" timestamp (row 16 / A16) refreshed to 10:22 AM

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H12 and Q12 are already formatted as Text (@), so plain assignment keeps
# them as text.
$ws.Range("H12").Value = "1:0"
$ws.Range("Q12").Value = "2:0"

# P12 is formatted with a numeric format (0.00) even though the workbook
# stores its value as literal text ("36.0000", 4 decimal places). Flip the
# cell to Text just long enough to assign the value verbatim, then restore
# its original number format so style/appearance stay unchanged.
$p12 = $ws.Range("P12")
$origFormat = $p12.NumberFormat
$p12.NumberFormat = "@"
$p12.Value = "36.0000"
$p12.NumberFormat = $origFormat

# Recalculated total for the sale-price column.
$ws.Range("P15").Value = 268.89999999999998

# Refresh the report generation timestamp.
$ws.Range("A16").Value = "Saturday, 13 September, 2025 10:22 AM"
